# Update the build timestamp embedded in version strings throughout the
# workbook: "February 03 2026 17.29.55 EST" -> "February 03 2026 18.05.36 EST"

$wb = $excel.ActiveWorkbook

$oldStamp = "February 03 2026 17.29.55 EST"
$newStamp = "February 03 2026 18.05.36 EST"

# --- "About" sheet -------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")

$rA2 = $wsAbout.Range("A2")
$rA2.Value = $rA2.Value().Replace($oldStamp, $newStamp)

$rA6 = $wsAbout.Range("A6")
$rA6.Value = $rA6.Value().Replace($oldStamp, $newStamp)

# --- "Boundaries and methane sources" sheet -------------------------------
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

for ($row = 2; $row -le 8; $row++) {
    $cell = $wsData.Cells.Item($row, 19)  # column S = build_version
    $cell.Value = $cell.Value().Replace($oldStamp, $newStamp)
}
